# Update gh-pages output (commit 456a3b4):
#  - "广州·天玑NOW x 崩三Only" event was cancelled: append （取消） to its name
#    and mark its min-price cell as "不可售" (not for sale) instead of a number.
#  - A handful of "想去人数" (interest-count) cells ticked up slightly.
#
# These edits touch three worksheets:
#   展览   (sheet index 1) - exhibitions
#   本地生活 (sheet index 3) - local-life events
#   全部类型 (sheet index 4) - the combined "all types" roll-up of every sheet
#
# 演出 (sheet index 2) has no changes in this revision.

$wb = $excel.ActiveWorkbook

$cancelledName = "广州·天玑NOW x 崩三Only（取消）"
$notForSale = "不可售"

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 3: 广州·天玑NOW x 崩三Only -> cancelled
$ws1.Range("C3").Value = $cancelledName
$ws1.Range("G3").Value = $notForSale

# Interest-count ("想去人数", column F) bumps
$ws1.Range("F6").Value = 370
$ws1.Range("F11").Value = 1350
$ws1.Range("F12").Value = 3008
$ws1.Range("F13").Value = 435
$ws1.Range("F14").Value = 1630
$ws1.Range("F15").Value = 1365
$ws1.Range("F16").Value = 799
$ws1.Range("F18").Value = 1395
$ws1.Range("F21").Value = 1128
$ws1.Range("F22").Value = 403
$ws1.Range("F23").Value = 4
$ws1.Range("F24").Value = 3498
$ws1.Range("F25").Value = 693

# ---------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Row 3: 广州·一人之下快闪店
$ws3.Range("F3").Value = 10

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (combined roll-up of every other sheet)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Row 5: 广州·一人之下快闪店
$ws4.Range("F5").Value = 10

# Row 6: 广州·天玑NOW x 崩三Only -> cancelled
$ws4.Range("C6").Value = $cancelledName
$ws4.Range("G6").Value = $notForSale

# Interest-count ("想去人数", column F) bumps
$ws4.Range("F16").Value = 370
$ws4.Range("F21").Value = 1350
$ws4.Range("F22").Value = 3008
$ws4.Range("F23").Value = 435
$ws4.Range("F24").Value = 1630
$ws4.Range("F25").Value = 1365
$ws4.Range("F26").Value = 799
$ws4.Range("F28").Value = 1395
$ws4.Range("F33").Value = 1128
$ws4.Range("F34").Value = 403
$ws4.Range("F35").Value = 4
$ws4.Range("F36").Value = 3498
$ws4.Range("F37").Value = 693
